$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct previously-reported daily counts (col B) for 12/8-12/17/2020 ---
$ws.Cells.Item(275, 2).Value = 41
$ws.Cells.Item(280, 2).Value = 53
$ws.Cells.Item(281, 2).Value = 51
$ws.Cells.Item(283, 2).Value = 49
$ws.Cells.Item(284, 2).Value = 52

# --- Append new row for 12/18/2020 ---
$ws.Cells.Item(285, 1).Value = 44183
$ws.Cells.Item(285, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(285, 2).Value = 42
$ws.Cells.Item(285, 3).Formula = "=B285+C284"
$ws.Cells.Item(285, 4).Value = 1
$ws.Cells.Item(285, 5).Formula = "=D285+E284"
$ws.Cells.Item(285, 6).Formula = "=AVERAGE(B279:B285)"

# --- Append new row for 12/19/2020 ---
$ws.Cells.Item(286, 1).Value = 44184
$ws.Cells.Item(286, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(286, 2).Value = 26
$ws.Cells.Item(286, 3).Formula = "=B286+C285"
$ws.Cells.Item(286, 4).Value = 0
$ws.Cells.Item(286, 5).Formula = "=D286+E285"
$ws.Cells.Item(286, 6).Formula = "=AVERAGE(B280:B286)"

# --- View state: keep header row frozen, scroll near the new bottom rows,
#     and leave the active selection on K284 like the saved workbook ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 263
$win.ScrollColumn = 2
$ws.Range("K284").Select()
